# Add new Testcase for OneWayBooking Full
# Adds a new worksheet "OneWayBooking" after "LoginData" and populates it
# with the one-way booking test fixture rows.

$wb = $excel.ActiveWorkbook
$loginData = $wb.Worksheets.Item(1)

# Insert the new sheet right after LoginData; Worksheets.Add() makes it the
# active sheet, which matches the target's activeTab/tabSelected change.
$ws = $wb.Worksheets.Add($null, $loginData)
$ws.Name = "OneWayBooking"

# ---- Header / data rows ------------------------------------------------
$rows = @(
    @("agileway", "testwise", "oneway", "Joseph", "Lieven",   "visa",   4659908765341267),
    @("agileway", "testwise", "oneway", "Chris",  "Lieven",   "master", 4659908765341267),
    @("agileway", "testwise", "oneway", "Joseph", "Hamilton", "visa",   4659908765341267),
    @("agileway", "testwise", "oneway", "John",   "Dobbson",  "master", 4659908765341267),
    @("agileway", "testwise", "oneway", "Joseph", "Wilson",   "visa",   4659908765341267)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowNum = $r + 1
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$c]
    }
    $ws.Rows.Item($rowNum).RowHeight = 13.8
}

# ---- Column G formatting / validation ----------------------------------
$ws.Columns.Item(7).ColumnWidth = 24.6

# G2: custom / between validation (added first to match sqref ordering of
# the committed workbook's dataValidations list: G2 then G1)
$ws.Range("G2").Validation.Add(7, 1, 1, "0", "0")
$ws.Range("G2").Validation.IgnoreBlank = $true
$ws.Range("G2").Validation.ShowInput = $false
$ws.Range("G2").Validation.ShowError = $true

# G1: text-length validation
$ws.Range("G1").Validation.Add(6, 1, 3, "0", "0")
$ws.Range("G1").Validation.IgnoreBlank = $true
$ws.Range("G1").Validation.ShowInput = $false
$ws.Range("G1").Validation.ShowError = $true

# ---- Page setup ----------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Zoom = 100
$ws.PageSetup.Order = 1
$ws.PageSetup.Orientation = 1
$ws.PageSetup.BlackAndWhite = $false
$ws.PageSetup.Draft = $false
$ws.PageSetup.PrintComments = 0
$ws.PageSetup.PrintHeadings = $false
$ws.PageSetup.PrintGridlines = $false
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&Kffffff&A"
$ws.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12&KffffffPage &P"

# ---- Selection state matching the committed workbook --------------------
$ws.Range("I9").Select()
